# EPBDS-2847 Generate custom typed SpreadsheetResult.
# Added possibility to test spreadsheet cells via test table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets 1 & 2 (Cyrillic placeholders -> meaningful test names)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "Test Datatype Array"
$ws2.Name = "Test Array"

Write-Output "renamed sheets"

# ---------------------------------------------------------------------------
# 2. Build "Test Array" (sheet2) content: a Spreadsheet test-table that
#    exercises a custom typed SpreadsheetResult (DoubleValue[] / AnyValue).
#    The layout mirrors the existing "Test Datatype Array" sheet, so we copy
#    the matching cell *formats* from there and just change the text.
# ---------------------------------------------------------------------------

# Column widths / row layout
$ws2.Columns.Item(3).ColumnWidth = 16.85546875
$ws2.Columns.Item(5).ColumnWidth = 34
$ws2.Columns.Item(9).ColumnWidth = 70.140625

# Row 5 : title bar (merged C5:E5) + method header in I5
$ws1.Range("C7:E7").Copy()
$ws2.Range("C5:E5").PasteSpecial(-4122)
$ws2.Range("C5:E5").Merge()
$ws2.Range("C5").Value = "Spreadsheet DoubleValue[] testDoubleValue()`n"
$ws2.Rows.Item(5).RowHeight = 25.5

$ws1.Range("C10").Copy()
$ws2.Range("I5").PasteSpecial(-4122)
$ws2.Range("I5").Value = "Method DoubleValue[] getDVs()"

Write-Output "row5 done"

